$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_specific_translations")
$ws.Activate()

# Fix the Portuguese translation of the HHID prompt (was left in English)
$ws.Range("C3").Value = "Identificação do agregado: {{data.hh_id}}"

# Add the missing Kiswahili translation of the HHID prompt
$ws.Range("D3").Value = "Utambulisho wa Kaya: {{data.hh_id}}"

$ws.Range("C14").Select()
